# blind_75.xlsx edit: rename "Notes" columns to "Approach" columns (BF + MO),
# move the BF rich-text note into the new "BF Approach" cell, and fill the
# "BF Solution" cell with a link to the solution source file on GitHub.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Stash the BF Notes rich text (H2) in a scratch cell far away so we
#        can move it (with its bold formatting intact) into the new BF
#        Approach cell (F2) later, after we've shifted the other BF cells
#        right to make room. ---
$ws.Range("H2").Copy()
$ws.Range("Z1").PasteSpecial()
$ws.Range("H2").Copy()
$ws.Range("Z1").PasteSpecial()

# --- 2. Shift the BF block one column to the right (F,G -> G,H), freeing up
#        column F for the new "BF Approach" column. Go right-to-left so we
#        don't clobber values before they're read. ---
$ws.Range("H1").Value2 = $ws.Range("G1").Value2
$ws.Range("H2").Value2 = $ws.Range("G2").Value2

$ws.Range("G1").Value2 = $ws.Range("F1").Value2
$ws.Range("G2").Value2 = $ws.Range("F2").Value2

$ws.Range("F1").Value2 = "BF Approach"
$ws.Range("F2").ClearContents()
$ws.Range("Z1").Copy()
$ws.Range("F2").PasteSpecial()
$ws.Range("Z1").Clear()

# --- 3. BF Solution (E2) now holds a link to the brute-force solution file. ---
$ws.Range("E2").Value2 = "https://github.com/davidmeadejr/leetcode/blob/main/blind_75/array/brute_force/two_sum.py"

# --- 4. Shift the MO block one column to the right (J,K -> K,L), freeing up
#        column J for the new "MO Approach" column. ---
$ws.Range("L1").Value2 = $ws.Range("K1").Value2
$ws.Range("K1").Value2 = $ws.Range("J1").Value2
$ws.Range("J1").Value2 = "MO Approach"

# --- 5. Row-7 leftover formatting markers: drop the one in D7, and shift the
#        one in G7 to H7 (it moved with the rest of column G -> H). ---
$ws.Range("H7").WrapText = $true
$ws.Range("D7").Clear()
$ws.Range("G7").Clear()

# --- 6. Styling touch-ups. ---
$ws.Range("K2").WrapText = $true
$ws.Range("K2").VerticalAlignment = -4108

$ws.Rows(2).RowHeight = 102

$ws.Columns("F").ColumnWidth = 21.17
$ws.Columns("J").ColumnWidth = 12.33

# --- 7. Leave the cursor where the author left it. ---
$ws.Range("E2").Select()

Write-Host "edit complete"
